# Add a "controller vibration" option row to the Strings sheet.
#
# A new row is inserted immediately below the "Invert look Y" row (row 60,
# pushing "look sensitivity" and everything after it down by one row), with
# the new localization key/value pair for the controller vibration setting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60 (shifts existing row 60 "look sensitivity" and
# all subsequent rows down by one). The inserted row inherits formatting from
# the row above it (row 59), matching the plain/unwrapped style already used
# throughout this column.
$ws.Rows.Item(60).Insert()

# Populate the new row with the controller-vibration localization strings.
$ws.Cells.Item(60, 1).Value = "controller vibration"
$ws.Cells.Item(60, 2).Value = "Controller vibration"

# Match the row height used by every other (non-wrapped) row in this table.
$ws.Rows.Item(60).RowHeight = 13.4

# Update the saved view state to match (scrolled/selected position after the edit).
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A61").Select()
